# Checklist PvA: mark additional items as checked, clear a couple of
# stray checkmarks/question-marks that were leftovers from a manual
# review round, and drop the "Controle door groep" remarks now that
# those items have been resolved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Row 13: was marked with a "?" in B13 and had "  " in D13 -> now
# checked (1) and annotated with "v" instead.
$ws.Range("B13").Value = 1
$ws.Range("D13").Value = "v"

# Row 40: was marked with a "v" text in B40 -> replace with a numeric
# check (1), same as all the other checked rows.
$ws.Range("B40").Value = 1

# Rows that get checked off (B = 1) and lose their
# "Controle door groep" remark in column D.
$rowsCheckedAndCleared = 26,27,28,29,31,49,51,72,75,77,78,85,86,93,94

foreach ($r in $rowsCheckedAndCleared) {
    $ws.Range("B$r").Value = 1
    $ws.Range("D$r").ClearContents()
}

# Row 100 only loses the remark in D; B100 was already checked.
$ws.Range("D100").ClearContents()

# Row 32: this item goes back to being unchecked.
$ws.Range("B32").ClearContents()

$wb.Save()
